$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("AK_1.png", 0.502, 0.002, 0.982, 0.333, "Tidak Diketahui", "Salah"),
    @("AK_2.png", 0.883, 0.003, 0.988, 0.4, "Tidak Diketahui", "Salah"),
    @("AK_3.png", 0.853, 0.003, 0.979, 0.8, "Akhlak Kamiswara", "Benar"),
    @("AK_4.png", 0.533, 0.002, 0.976, 0.533, "Akhlak Kamiswara", "Benar"),
    @("AK_5.png", 0.531, 0.002, 0.99, 0.333, "Tidak Diketahui", "Salah"),
    @("MIB_1.png", 1.553, 0.005, 0.976, 0.667, "Muhammad Iqbal Baqi", "Benar"),
    @("MIB_2.png", 1.159, 0.004, 0.973, 0.467, "Tidak Diketahui", "Salah"),
    @("MIB_3.png", 1.698, 0.006, 0.999, 1, "Muhammad Iqbal Baqi", "Benar"),
    @("MIB_4.png", 1.174, 0.004, 0.985, 0.733, "Muhammad Iqbal Baqi", "Benar"),
    @("MIB_5.png", 1.495, 0.005, 0.983, 0.6, "Muhammad Iqbal Baqi", "Benar"),
    @("AAH_1.png", 0.872, 0.003, 0.975, 0.533, "Andrea Ayunove Hutami", "Benar"),
    @("AAH_2.png", 1.103, 0.004, 0.992, 0.6, "Andrea Ayunove Hutami", "Benar"),
    @("AAH_3.png", 0.907, 0.003, 0.973, 0.667, "Andrea Ayunove Hutami", "Benar"),
    @("TI_1.png", 0.9379999999999999, 0.003, 0.949, 0.533, "Toni Ismail", "Benar"),
    @("TI_2.png", 0.9419999999999999, 0.003, 0.989, 0.333, "Tidak Diketahui", "Salah"),
    @("TI_3.png", 0.752, 0.003, 0.991, 0.6, "Toni Ismail", "Benar"),
    @("TI_4.png", 0.746, 0.002, 0.974, 0.733, "Toni Ismail", "Benar"),
    @("TI_5.png", 1.04, 0.003, 0.992, 0.467, "Tidak Diketahui", "Salah"),
    @("RAS_1.png", 0.675, 0.002, 0.975, 0.333, "Tidak Diketahui", "Salah"),
    @("RAS_2.png", 1.007, 0.003, 0.966, 0.6, "Ridha Ayu Salsabila", "Benar"),
    @("RAS_3.png", 0.57, 0.002, 0.982, 0.4, "Tidak Diketahui", "Salah"),
    @("RAS_4.png", 1.259, 0.004, 0.975, 0.267, "Tidak Diketahui", "Salah"),
    @("RAS_5.png", 1.186, 0.004, 0.971, 0.8, "Ridha Ayu Salsabila", "Benar"),
    @("RR_1.png", 1.441, 0.005, 0.981, 0.467, "Tidak Diketahui", "Salah"),
    @("RR_2.png", 1.361, 0.005, 0.987, 0.333, "Tidak Diketahui", "Salah"),
    @("RR_3.png", 1.103, 0.004, 0.906, 0.4, "Tidak Diketahui", "Salah"),
    @("RR_4.png", 1.376, 0.005, 0.989, 0.4, "Tidak Diketahui", "Salah"),
    @("RR_5.png", 1.417, 0.005, 0.984, 0.667, "Rafiqo Rapitasari", "Benar"),
    @("AR_1.png", 0.891, 0.003, 0.965, 0.733, "Arizli Romadhon", "Benar"),
    @("GA_1.png", 1.548, 0.005, 0.983, 0.9330000000000001, "Gege Ardiyansyah", "Benar"),
    @("GA_2.png", 0.676, 0.002, 0.99, 0.267, "Tidak Diketahui", "Salah"),
    @("GA_3.png", 0.721, 0.002, 0.99, 0.333, "Tidak Diketahui", "Salah"),
    @("FY_1.png", 1.11, 0.004, 0.969, 0.333, "Tidak Diketahui", "Salah"),
    @("FY_2.png", 1.589, 0.005, 0.979, 0.333, "Tidak Diketahui", "Salah"),
    @("FY_3.png", 1.452, 0.005, 0.991, 0.4, "Tidak Diketahui", "Salah"),
    @("FY_4.png", 1.227, 0.004, 0.984, 0.267, "Tidak Diketahui", "Salah"),
    @("TO_1.png", 0.806, 0.003, 0.971, 0.4, "Tidak Diketahui", "Salah"),
    @("TO_2.png", 1.065, 0.004, 0.986, 0.533, "Tiara Oktavian", "Benar"),
    @("TO_3.png", 0.834, 0.003, 0.983, 0.6, "Tiara Oktavian", "Benar"),
    @("TO_4.png", 3.505, 0.011, 0.841, 0.4, "Tidak Diketahui", "Salah"),
    @("TO_5.png", 3.187, 0.01, 0.841, 0.467, "Tidak Diketahui", "Salah"),
    @("TD_1.png", 2.35, 0.008, 0.858, 0.733, "Muhammad Iqbal Baqi", "Salah"),
    @("TD_2.png", 2.49, 0.008, 0.854, 0.4, "Tidak Diketahui", "Benar"),
    @("TD_3.png", 1.058, 0.004, 0.955, 0.267, "Tidak Diketahui", "Benar"),
    @("TD_4.png", 0.998, 0.003, 0.987, 0.267, "Tidak Diketahui", "Benar"),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    for ($j = 0; $j -lt $data[$i].Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}

Write-Host "Done writing data"